$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update header values for columns B-E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: D2 value removed; new values placed in B2/C2
$ws.Range("D2").ClearContents()
$ws.Range("B2").Value = 30.373805491377226
$ws.Range("C2").Value = 27.783653487400766

# Row 3: B3 value removed; C3 value updated
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 25.230453533488497

# Update the active selection to match the new data extent
$ws.Range("B1:E3").Select()
